$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure rows 7-9 have the same style as existing data rows (copy format from row 2).
# Column C is intentionally skipped (it has no data in any row, and copying its
# blank format would create stray empty <c> cell records).
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A7:B9").PasteSpecial(-4122) | Out-Null
$ws.Range("D2:W2").Copy() | Out-Null
$ws.Range("D7:W9").PasteSpecial(-4122) | Out-Null

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("D2").Value = 228.5
$ws.Range("E2").Value = 19322371
$ws.Range("F2").Value = 113.3658536585366
$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = "Washington"
$ws.Range("I2").Value = "Chicago"
$ws.Range("J2").Value = 0.4939024390243902
$ws.Range("K2").Value = 98.75243902439024
$ws.Range("L2").Value = 113.9414634146341
$ws.Range("M2").Value = 114.9731707317073
$ws.Range("N2").Value = 77.76707317073171
$ws.Range("O2").Value = 0.3463536585365854
$ws.Range("P2").Value = 0.5854634146341462
$ws.Range("Q2").Value = 0.2659634146341464
$ws.Range("R2").Value = 12.09390243902439
$ws.Range("S2").Value = 11.74024390243903
$ws.Range("T2").Value = 0.2083841463414634
$ws.Range("U2").Value = 0.996185005786789
$ws.Range("V2").Value = 1.008397028021033
$ws.Range("W2").Value = 10.44546646695767

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("D3").Value = 232
$ws.Range("E3").Value = 5116589
$ws.Range("F3").Value = 113.2619179600887
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = "Detroit"
$ws.Range("I3").Value = "Minnesota"
$ws.Range("J3").Value = 0.5340328984685195
$ws.Range("K3").Value = 100.0564024390244
$ws.Range("L3").Value = 112.7019955654102
$ws.Range("M3").Value = 116.6864745011086
$ws.Range("N3").Value = 74.14212860310421
$ws.Range("O3").Value = 0.3745354767184035
$ws.Range("P3").Value = 0.57735088691796
$ws.Range("Q3").Value = 0.2991524390243901
$ws.Range("R3").Value = 13.11798780487805
$ws.Range("S3").Value = 12.46804323725055
$ws.Range("T3").Value = 0.2282720343680709
$ws.Range("U3").Value = 0.9952716868197601
$ws.Range("V3").Value = 1.002093628025482
$ws.Range("W3").Value = 11.26227294766619

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("D4").Value = 228.5
$ws.Range("E4").Value = 7962208.5
$ws.Range("F4").Value = 113.8125
$ws.Range("G4").Value = 4.5
$ws.Range("H4").Value = "Atlanta"
$ws.Range("I4").Value = "Milwaukee"
$ws.Range("J4").Value = 0.55
$ws.Range("K4").Value = 99.68500000000002
$ws.Range("L4").Value = 112.94375
$ws.Range("M4").Value = 112.8325
$ws.Range("N4").Value = 76.79749999999999
$ws.Range("O4").Value = 0.381875
$ws.Range("P4").Value = 0.5632999999999999
$ws.Range("Q4").Value = 0.25385
$ws.Range("R4").Value = 11.83375
$ws.Range("S4").Value = 11.53125
$ws.Range("T4").Value = 0.19944375
$ws.Range("U4").Value = 1.000109841827768
$ws.Range("V4").Value = 0.9808954594845883
$ws.Range("W4").Value = 10.6887837355695

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("D5").Value = 231
$ws.Range("E5").Value = 12984634.5
$ws.Range("F5").Value = 117.9634146341463
$ws.Range("G5").Value = 9.5
$ws.Range("H5").Value = "Boston"
$ws.Range("I5").Value = "NewOrleans"
$ws.Range("J5").Value = 0.5490931832395247
$ws.Range("K5").Value = 99.15731707317073
$ws.Range("L5").Value = 117.6634146341463
$ws.Range("M5").Value = 112.4219512195122
$ws.Range("N5").Value = 77.84146341463415
$ws.Range("O5").Value = 0.410890243902439
$ws.Range("P5").Value = 0.5952073170731706
$ws.Range("Q5").Value = 0.2806463414634147
$ws.Range("R5").Value = 12.13292682926829
$ws.Range("S5").Value = 12.44390243902439
$ws.Range("T5").Value = 0.2097560975609756
$ws.Range("U5").Value = 1.036585365853659
$ws.Range("V5").Value = 0.9932591518117657
$ws.Range("W5").Value = 10.05097842903586

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("D6").Value = 226
$ws.Range("E6").Value = 9262082.5
$ws.Range("F6").Value = 114.8048780487805
$ws.Range("G6").Value = 4.5
$ws.Range("H6").Value = "NewYork"
$ws.Range("I6").Value = "Indiana"
$ws.Range("J6").Value = 0.4875
$ws.Range("K6").Value = 98.6731707317073
$ws.Range("L6").Value = 115.620731707317
$ws.Range("M6").Value = 114.1219512195122
$ws.Range("N6").Value = 74.58902439024391
$ws.Range("O6").Value = 0.4136585365853658
$ws.Range("P6").Value = 0.5721585365853659
$ws.Range("Q6").Value = 0.2745487804878048
$ws.Range("R6").Value = 11.76219512195122
$ws.Range("S6").Value = 12.29634146341463
$ws.Range("T6").Value = 0.2201646341463414
$ws.Range("U6").Value = 1.008830211324961
$ws.Range("V6").Value = 0.9957977902620151
$ws.Range("W6").Value = 11.98146739930713

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2023
$ws.Range("D7").Value = 238.5
$ws.Range("E7").Value = 9992204.5
$ws.Range("F7").Value = 114.1993902439024
$ws.Range("G7").Value = 13.5
$ws.Range("H7").Value = "Memphis"
$ws.Range("I7").Value = "SanAntonio"
$ws.Range("J7").Value = 0.5317139001349528
$ws.Range("K7").Value = 100.4958231707317
$ws.Range("L7").Value = 113.29875
$ws.Range("M7").Value = 115.0175
$ws.Range("N7").Value = 75.95158536585365
$ws.Range("O7").Value = 0.353910975609756
$ws.Range("P7").Value = 0.5622545731707318
$ws.Range("Q7").Value = 0.2565567073170731
$ws.Range("R7").Value = 12.55484756097561
$ws.Range("S7").Value = 12.25079268292683
$ws.Range("T7").Value = 0.1932905487804878
$ws.Range("U7").Value = 1.003509580350636
$ws.Range("V7").Value = 1.046319236921261
$ws.Range("W7").Value = 11.35277984872319

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2023
$ws.Range("D8").Value = 226.5
$ws.Range("E8").Value = 6049076.5
$ws.Range("F8").Value = 114.9244047619048
$ws.Range("G8").Value = 14
$ws.Range("H8").Value = "Denver"
$ws.Range("I8").Value = "Phoenix"
$ws.Range("J8").Value = 0.4817073170731707
$ws.Range("K8").Value = 97.29541666666665
$ws.Range("L8").Value = 117.3911904761905
$ws.Range("M8").Value = 114.871369047619
$ws.Range("N8").Value = 76.4417261904762
$ws.Range("O8").Value = 0.3649839285714286
$ws.Range("P8").Value = 0.590007142857143
$ws.Range("Q8").Value = 0.2486589285714286
$ws.Range("R8").Value = 12.44607142857143
$ws.Range("S8").Value = 12.24357142857143
$ws.Range("T8").Value = 0.2096875000000001
$ws.Range("U8").Value = 1.009880533935894
$ws.Range("V8").Value = 0.991124236836904
$ws.Range("W8").Value = 11.34684067436913

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 2023
$ws.Range("D9").Value = 236.5
$ws.Range("E9").Value = 18084349.5
$ws.Range("F9").Value = 113.9349358974359
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = "Sacramento"
$ws.Range("I9").Value = "Houston"
$ws.Range("J9").Value = 0.4743421052631579
$ws.Range("K9").Value = 99.6941346153846
$ws.Range("L9").Value = 113.7744871794872
$ws.Range("M9").Value = 116.4991666666667
$ws.Range("N9").Value = 76.89932692307693
$ws.Range("O9").Value = 0.4110048076923076
$ws.Range("P9").Value = 0.5767176282051283
$ws.Range("Q9").Value = 0.2881881410256409
$ws.Range("R9").Value = 13.24182692307692
$ws.Range("S9").Value = 12.29208333333333
$ws.Range("T9").Value = 0.2205663461538462
$ws.Range("U9").Value = 1.00118572844847
$ws.Range("V9").Value = 1.028856924773938
$ws.Range("W9").Value = 11.25332596191566

